$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing scalar odds values (row 2) ---
$ws.Range("F2").Value = 2.88
$ws.Range("N2").Value = 3.55
$ws.Range("P2").Value = 1.88
$ws.Range("R2").Value = 1.34
$ws.Range("U2").Value = 2.12
$ws.Range("AC2").Value = 8.199999999999999

# --- Update existing scalar odds values (row 3) ---
$ws.Range("Y3").Value = 22
$ws.Range("AN3").Value = 8.4

# --- Update existing scalar odds values (row 6) ---
$ws.Range("S6").Value = 2.2
$ws.Range("T6").Value = 1.78

# --- Update existing scalar odds values (row 8) ---
$ws.Range("F8").Value = 1.71
$ws.Range("G8").Value = 1.72
$ws.Range("I8").Value = 6.2
$ws.Range("P8").Value = 1.92
$ws.Range("T8").Value = 2.02
$ws.Range("V8").Value = 1.19
$ws.Range("W8").Value = 2.38
$ws.Range("AC8").Value = 8.6
$ws.Range("AJ8").Value = 16.5

# --- Append new row 13: Necaxa vs Monterrey ---
$ws.Range("A13").Value = "Mexican Liga MX"
# Date-looking text must stay text, not get auto-converted to a date serial.
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2026-01-13"
$ws.Range("B13").ClearFormats()
$ws.Range("C13").Value = "22:00:00"
$ws.Range("D13").Value = "Necaxa"
$ws.Range("E13").Value = "Monterrey"
$ws.Range("F13").Value = 2.58
$ws.Range("G13").Value = 2.74
$ws.Range("H13").Value = 2.64
$ws.Range("I13").Value = 2.8
$ws.Range("J13").Value = 3.8
$ws.Range("K13").Value = 4
$ws.Range("L13").Value = 1.01
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 1.21
$ws.Range("P13").Value = 2.38
$ws.Range("Q13").Value = 1.61
$ws.Range("R13").Value = 1.55
$ws.Range("S13").Value = 2.52
$ws.Range("T13").Value = 1.54
$ws.Range("U13").Value = 2.54
$ws.Range("V13").Value = 1.55
$ws.Range("W13").Value = 1.57
$ws.Range("X13").Value = 24
$ws.Range("Y13").Value = 16.5
$ws.Range("Z13").Value = 22
$ws.Range("AA13").Value = 40
$ws.Range("AB13").Value = 16
$ws.Range("AC13").Value = 10
$ws.Range("AD13").Value = 13.5
$ws.Range("AE13").Value = 27
$ws.Range("AF13").Value = 22
$ws.Range("AG13").Value = 13.5
$ws.Range("AH13").Value = 16
$ws.Range("AI13").Value = 34
$ws.Range("AJ13").Value = 40
$ws.Range("AK13").Value = 27
$ws.Range("AL13").Value = 34
$ws.Range("AM13").Value = 65
$ws.Range("AN13").Value = 16.5
$ws.Range("AO13").Value = 17

# --- Append new row 14: Pachuca vs Leon ---
$ws.Range("A14").Value = "Mexican Liga MX"
# Date-looking text must stay text, not get auto-converted to a date serial.
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "2026-01-13"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").Value = "22:06:00"
$ws.Range("D14").Value = "Pachuca"
$ws.Range("E14").Value = "Leon"
$ws.Range("F14").Value = 1.67
$ws.Range("G14").Value = 1.75
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = 5.1
$ws.Range("J14").Value = 4.3
$ws.Range("K14").Value = 4.8
$ws.Range("L14").Value = 1.01
$ws.Range("M14").Value = 1.01
$ws.Range("N14").Value = 2.32
$ws.Range("O14").Value = 1.21
$ws.Range("P14").Value = 2.32
$ws.Range("Q14").Value = 1.64
$ws.Range("R14").Value = 1.45
$ws.Range("S14").Value = 2.34
$ws.Range("T14").Value = 1.6
$ws.Range("U14").Value = 1.01
$ws.Range("V14").Value = 1.24
$ws.Range("W14").Value = 2.32
$ws.Range("X14").Value = 26
$ws.Range("Y14").Value = 30
$ws.Range("Z14").Value = 55
$ws.Range("AA14").Value = 1000
$ws.Range("AB14").Value = 13
$ws.Range("AC14").Value = 14.5
$ws.Range("AD14").Value = 29
$ws.Range("AE14").Value = 75
$ws.Range("AF14").Value = 14.5
$ws.Range("AG14").Value = 12.5
$ws.Range("AH14").Value = 22
$ws.Range("AI14").Value = 75
$ws.Range("AJ14").Value = 24
$ws.Range("AK14").Value = 19.5
$ws.Range("AL14").Value = 36
$ws.Range("AM14").Value = 1000
$ws.Range("AN14").Value = 1000
$ws.Range("AO14").Value = 1000
